$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the datetime value in A2 (bug fix: timestamp corrected)
$ws.Range("A2").Value = "10.05.2024 14:10:46"

# Remove the now-duplicate row 3 entirely (shrinks used range to A1:B2)
$ws.Rows(3).Delete()
